# 2nd submission. Revised based on comments.
#
# - Rename "Sheet1 (2)" -> "Model Refinements"
# - Keep "Sheet1" named "Sheet1" (its sheetId/rId shift because a brand-new
#   "Confusion Mat" sheet is inserted between the two existing sheets)
# - Insert a new "Confusion Mat" worksheet (with the cross-validation /
#   confusion-matrix tables) right after "Model Refinements"
# - Update a couple of remembered selections

$wb = $excel.ActiveWorkbook

$modelRefinements = $wb.Worksheets.Item(1)

# New sheet, inserted right after "Model Refinements" -> becomes tab #2.
# NOTE: grab this *before* renaming "Model Refinements" (Item(1) is a
# positional handle; renaming doesn't move it, but we keep the rename for
# last anyway to keep this section easy to reason about).
$confusionMat = $wb.Worksheets.Add($null, $modelRefinements)
$confusionMat.Name = "Confusion Mat"

$modelRefinements.Name = "Model Refinements"

# IMPORTANT: the original "Sheet1" was tab #2 before the insert, but
# inserting "Confusion Mat" right after "Model Refinements" pushes it to
# tab #3. A positional handle grabbed *before* the insert re-resolves to
# whatever now sits at that position (i.e. "Confusion Mat"), so fetch the
# "Sheet1" handle by name, after the insert, instead.
$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet1.Name = "Sheet1"

# ---------------------------------------------------------------------
# Confusion Mat content
# ---------------------------------------------------------------------

function Set-Cell($ws, $ref, $val) {
    $ws.Range($ref).Value = $val
}

# Header row (slot 1..5, AVG)
Set-Cell $confusionMat "C3" "slot 1"
Set-Cell $confusionMat "D3" "slot 2"
Set-Cell $confusionMat "E3" "slot 3"
Set-Cell $confusionMat "F3" "slot 4"
Set-Cell $confusionMat "G3" "slot 5"
Set-Cell $confusionMat "H3" "AVG"

# TP row
Set-Cell $confusionMat "B4" "TP"
Set-Cell $confusionMat "C4" 12066
Set-Cell $confusionMat "D4" 11951
Set-Cell $confusionMat "E4" 12617
Set-Cell $confusionMat "F4" 13022
Set-Cell $confusionMat "G4" 13066
$confusionMat.Range("H4").Formula = "=SUM(C4:G4)/5"

# FP row
Set-Cell $confusionMat "B5" "FP"
Set-Cell $confusionMat "C5" 1002
Set-Cell $confusionMat "D5" 1117
Set-Cell $confusionMat "E5" 451
Set-Cell $confusionMat "F5" 46
Set-Cell $confusionMat "G5" 2

# FN row
Set-Cell $confusionMat "B6" "FN"
Set-Cell $confusionMat "C6" 1002
Set-Cell $confusionMat "D6" 1857
Set-Cell $confusionMat "E6" 451
Set-Cell $confusionMat "F6" 46
Set-Cell $confusionMat "G6" 2

# Shared AVG formula column for FP/FN/Precision/Recall/F1 rows
$confusionMat.Range("H5:H9").Formula = "=SUM(C5:G5)/5"

# Precision row
Set-Cell $confusionMat "B7" "Precision"
$confusionMat.Range("C7").Formula = "=C4/(C4+C5)"
$confusionMat.Range("D7:G7").Formula = "=D4/(D4+D5)"

# Recall row
Set-Cell $confusionMat "B8" "Recall"
$confusionMat.Range("C8").Formula = "=C4/(C4+C6)"
$confusionMat.Range("D8:G8").Formula = "=D4/(D4+D6)"

# F1 row
Set-Cell $confusionMat "B9" "F1"
$confusionMat.Range("C9").Formula = "=2*(C7*C8)/(C7+C8)"
$confusionMat.Range("D9:G9").Formula = "=2*(D7*D8)/(D7+D8)"

# Precision/Recall/F1 are formatted with 4 decimal places
$confusionMat.Range("C7:H9").NumberFormat = "0.0000"

# ---------------------------------------------------------------------
# Per-slot confusion-matrix breakdown for FP ("FP" label, rows 11-22)
# ---------------------------------------------------------------------
Set-Cell $confusionMat "A11" "FP"

Set-Cell $confusionMat "B11" 0
Set-Cell $confusionMat "C11" 17
Set-Cell $confusionMat "D11" 122
Set-Cell $confusionMat "E11" 37
Set-Cell $confusionMat "F11" 5

Set-Cell $confusionMat "B12" 1
Set-Cell $confusionMat "C12" 214
Set-Cell $confusionMat "D12" 120
Set-Cell $confusionMat "E12" 63
Set-Cell $confusionMat "F12" 1

Set-Cell $confusionMat "B13" 2
Set-Cell $confusionMat "C13" 78
Set-Cell $confusionMat "D13" 63
Set-Cell $confusionMat "E13" 26
Set-Cell $confusionMat "F13" 3

Set-Cell $confusionMat "B14" 3
Set-Cell $confusionMat "C14" 128
Set-Cell $confusionMat "D14" 98
Set-Cell $confusionMat "E14" 37
Set-Cell $confusionMat "F14" 2

Set-Cell $confusionMat "B15" 4
Set-Cell $confusionMat "C15" 66
Set-Cell $confusionMat "D15" 51
Set-Cell $confusionMat "E15" 26
Set-Cell $confusionMat "F15" 3

Set-Cell $confusionMat "B16" 5
Set-Cell $confusionMat "C16" 102
Set-Cell $confusionMat "D16" 114
Set-Cell $confusionMat "E16" 46
Set-Cell $confusionMat "F16" 7

Set-Cell $confusionMat "B17" 6
Set-Cell $confusionMat "C17" 143
Set-Cell $confusionMat "D17" 103
Set-Cell $confusionMat "E17" 47
Set-Cell $confusionMat "F17" 1

Set-Cell $confusionMat "B18" 7
Set-Cell $confusionMat "C18" 62
Set-Cell $confusionMat "D18" 75
Set-Cell $confusionMat "E18" 20
Set-Cell $confusionMat "F18" 1

Set-Cell $confusionMat "B19" 8
Set-Cell $confusionMat "C19" 99
Set-Cell $confusionMat "D19" 106
Set-Cell $confusionMat "E19" 19
Set-Cell $confusionMat "F19" 3

Set-Cell $confusionMat "B20" 9
Set-Cell $confusionMat "C20" 93
Set-Cell $confusionMat "D20" 90
Set-Cell $confusionMat "E20" 39
Set-Cell $confusionMat "F20" 5

Set-Cell $confusionMat "B21" 10
Set-Cell $confusionMat "D21" 175
Set-Cell $confusionMat "E21" 91
Set-Cell $confusionMat "F21" 15
Set-Cell $confusionMat "G21" 2

$confusionMat.Range("C22").Formula = "=SUM(C11:C20)"
$confusionMat.Range("D22").Formula = "=SUM(D11:D21)"
$confusionMat.Range("E22").Formula = "=SUM(E11:E21)"
$confusionMat.Range("F22").Formula = "=SUM(F11:F21)"
$confusionMat.Range("G22").Formula = "=SUM(G11:G21)"

# ---------------------------------------------------------------------
# Per-slot confusion-matrix breakdown for FN ("FN" label, rows 24-35)
# ---------------------------------------------------------------------
Set-Cell $confusionMat "A24" "FN"

Set-Cell $confusionMat "B24" 0
Set-Cell $confusionMat "C24" 11
Set-Cell $confusionMat "D24" 94
Set-Cell $confusionMat "E24" 34
Set-Cell $confusionMat "F24" 5
Set-Cell $confusionMat "G24" 1

Set-Cell $confusionMat "B25" 1
Set-Cell $confusionMat "C25" 136
Set-Cell $confusionMat "D25" 88
Set-Cell $confusionMat "E25" 23
Set-Cell $confusionMat "F25" 4

Set-Cell $confusionMat "B26" 2
Set-Cell $confusionMat "C26" 205
Set-Cell $confusionMat "D26" 131
Set-Cell $confusionMat "E26" 32
Set-Cell $confusionMat "F26" 3

Set-Cell $confusionMat "B27" 3
Set-Cell $confusionMat "C27" 196
Set-Cell $confusionMat "D27" 151
Set-Cell $confusionMat "E27" 39
Set-Cell $confusionMat "F27" 5

Set-Cell $confusionMat "B28" 4
Set-Cell $confusionMat "C28" 79
Set-Cell $confusionMat "D28" 96
Set-Cell $confusionMat "E28" 22
Set-Cell $confusionMat "F28" 2

Set-Cell $confusionMat "B29" 5
Set-Cell $confusionMat "C29" 69
Set-Cell $confusionMat "D29" 96
Set-Cell $confusionMat "E29" 29
Set-Cell $confusionMat "F29" 2

Set-Cell $confusionMat "B30" 6
Set-Cell $confusionMat "C30" 79
Set-Cell $confusionMat "D30" 102
Set-Cell $confusionMat "E30" 22
Set-Cell $confusionMat "F30" 5

Set-Cell $confusionMat "B31" 7
Set-Cell $confusionMat "C31" 73
Set-Cell $confusionMat "D31" 60
Set-Cell $confusionMat "E31" 19
Set-Cell $confusionMat "F31" 1
Set-Cell $confusionMat "G31" 1

Set-Cell $confusionMat "B32" 8
Set-Cell $confusionMat "C32" 109
Set-Cell $confusionMat "D32" 140
Set-Cell $confusionMat "E32" 45
Set-Cell $confusionMat "F32" 3

Set-Cell $confusionMat "B33" 9
Set-Cell $confusionMat "C33" 45
Set-Cell $confusionMat "D33" 848
Set-Cell $confusionMat "E33" 32
Set-Cell $confusionMat "F33" 3

Set-Cell $confusionMat "B34" 10
Set-Cell $confusionMat "D34" 51
Set-Cell $confusionMat "E34" 154
Set-Cell $confusionMat "F34" 13

$confusionMat.Range("C35").Formula = "=SUM(C24:C33)"
$confusionMat.Range("D35").Formula = "=SUM(D24:D34)"
$confusionMat.Range("E35").Formula = "=SUM(E24:E34)"
$confusionMat.Range("F35").Formula = "=SUM(F24:F34)"
$confusionMat.Range("G35").Formula = "=SUM(G24:G34)"

# ---------------------------------------------------------------------
# Selections / active cells
# ---------------------------------------------------------------------

$confusionMat.Range("K16").Select()

$sheet1.Activate()
$sheet1.Range("F27").Select()

$modelRefinements.Activate()
$modelRefinements.Range("B57").Select()
